# Included mezz calculations + summary
# Updates the borrower cash-flow model on the "Cash Flows" sheet so the
# borrower_balance / margin_call / cash_balance / borrower_balance_change /
# borrower_cash_flow / borrower_cf_cum columns reflect the newly-included
# mezzanine debt draws (and their knock-on effects through the cumulative
# cash-flow summary).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 14317259.53416658
$ws.Range("F7").Value = -14317259.53416658
$ws.Range("N7").Value = -14317259.53416658
$ws.Range("Q7").Value = -20261143.63626546
$ws.Range("R7").Value = -48104907.58364561
$ws.Range("D8").Value = 29818032.94573001
$ws.Range("F8").Value = -15500773.41156343
$ws.Range("N8").Value = -15500773.41156343
$ws.Range("Q8").Value = -21543828.15932235
$ws.Range("R8").Value = -69648735.74296796
$ws.Range("D9").Value = 45761494.75871766
$ws.Range("F9").Value = -15943461.81298766
$ws.Range("N9").Value = -15943461.81298766
$ws.Range("Q9").Value = -22074757.07905615
$ws.Range("R9").Value = -91723492.82202412
$ws.Range("D10").Value = 60661770.50911177
$ws.Range("F10").Value = -14900275.75039411
$ws.Range("N10").Value = -14900275.75039411
$ws.Range("Q10").Value = -21108960.58486378
$ws.Range("R10").Value = -112832453.4068879
$ws.Range("D11").Value = 60661770.50911177
$ws.Range("R11").Value = -119107756.0372924
$ws.Range("D12").Value = 68539984.33019322
$ws.Range("N12").Value = -7878213.821081452
$ws.Range("Q12").Value = -14209441.65239663
$ws.Range("R12").Value = -133317197.6896891
$ws.Range("D13").Value = 76522795.20482963
$ws.Range("R13").Value = -147676548.1789692
$ws.Range("D14").Value = 83934451.55832386
$ws.Range("R14").Value = -161499521.6902956
$ws.Range("D15").Value = 91531813.19470388
$ws.Range("R15").Value = -175532522.9649983
$ws.Range("D16").Value = 98335814.46764934
$ws.Range("R16").Value = -188786110.4715011
$ws.Range("D17").Value = 98335814.46764934
$ws.Range("R17").Value = -195239346.5924793
$ws.Range("D18").Value = 98335814.46764934
$ws.Range("R18").Value = -201686015.0705067
$ws.Range("D19").Value = 98335814.46764934
$ws.Range("R19").Value = -208115977.5526537
$ws.Range("D20").Value = 98335814.46764934
$ws.Range("R20").Value = -214519174.8634328
$ws.Range("D21").Value = 98335814.46764934
$ws.Range("R21").Value = -220885627.0047986
$ws.Range("D22").Value = 98335814.46764934
$ws.Range("R22").Value = -227205433.1561478
$ws.Range("D23").Value = 98335814.46764934
$ws.Range("R23").Value = -233468771.674319
$ws.Range("D24").Value = 98335814.46764934
$ws.Range("R24").Value = -239665900.0935933
$ws.Range("D25").Value = 98335814.46764934
$ws.Range("R25").Value = -245586436.3029024
$ws.Range("D26").Value = 98335814.46764934
$ws.Range("R26").Value = -250223170.719888
$ws.Range("D27").Value = 98335814.46764934
$ws.Range("R27").Value = -254669497.364251
$ws.Range("D28").Value = 98335814.46764934
$ws.Range("R28").Value = -258938834.9622486
$ws.Range("D29").Value = 98335814.46764934
$ws.Range("R29").Value = -263031474.9463831
$ws.Range("D30").Value = 98335814.46764934
$ws.Range("R30").Value = -266960998.9564626
$ws.Range("D31").Value = 98335814.46764934
$ws.Range("R31").Value = -270725124.8308932
$ws.Range("D32").Value = 98335814.46764934
$ws.Range("R32").Value = -271533854.5380271
$ws.Range("D33").Value = 98335814.46764934
$ws.Range("R33").Value = -272458167.9299496
$ws.Range("D34").Value = 98335814.46764934
$ws.Range("R34").Value = -273507015.6979589
$ws.Range("D35").Value = 98335814.46764934
$ws.Range("R35").Value = -274689269.3559112
$ws.Range("D36").Value = 98335814.46764934
$ws.Range("R36").Value = -276013721.2402206
$ws.Range("D37").Value = 98335814.46764934
$ws.Range("R37").Value = -277489084.5098591
$ws.Range("D38").Value = 44816535.90076101
$ws.Range("N38").Value = 53519278.56688833
$ws.Range("Q38").Value = 51884369.93039092
$ws.Range("R38").Value = -225604714.5794682
$ws.Range("E39").Value = 61951990.08296807
$ws.Range("N39").Value = 44816535.90076101
$ws.Range("Q39").Value = 43013527.09331683
$ws.Range("E40").Value = 163164919.7441755
$ws.Range("E41").Value = 258523267.6063859
$ws.Range("R41").Value = -186735328.9430208
$ws.Range("E42").Value = 347729717.2807562
$ws.Range("E43").Value = 430488626.9302343
$ws.Range("R43").Value = -191652549.2966788
$ws.Range("E44").Value = 506506034.7516059
$ws.Range("R44").Value = -194421609.8097008
$ws.Range("E45").Value = 575489664.4754883
$ws.Range("R45").Value = -197408436.5919126
$ws.Range("E46").Value = 637148930.8843291
$ws.Range("R46").Value = -200621030.2053073
$ws.Range("E47").Value = 691194945.3484702
$ws.Range("R47").Value = -204067312.034436
$ws.Range("E48").Value = 737340521.3803352

Write-Output "applied 98 cell updates"